$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.723.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.75%  "
$ws.Range("D3").Value = "'1.729.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.23%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'227.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.40%  "
$ws.Range("D6").Value = "'0.5438"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.03%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'0.2729"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.53%  "
$ws.Range("D9").Value = "'0.06662"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.26%  "
$ws.Range("D10").Value = "'21.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.55%  "
$ws.Range("D11").Value = "'0.07756"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").Value = "'4.683"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("D13").Value = "'1.725.10"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.28%  "
$ws.Range("D14").Value = "'1.967.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.21%  "
$ws.Range("E15").Value = "  +5.21%  "
$ws.Range("D16").Value = "'0.0₅8378"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.23%  "
$ws.Range("D17").Value = "'68.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.58%  "
$ws.Range("D18").Value = "'27.710.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.76%  "
$ws.Range("D19").Value = "'225.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +17.45%  "
$ws.Range("D20").Value = "'4.799"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.21%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").Value = "'10.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E23").Value = "  +3.07%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").Value = "'146.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("D26").Value = "'1.726"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +13.06%  "
$ws.Range("E27").Value = "  +3.48%  "
$ws.Range("D28").Value = "'7.447"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.99%  "
$ws.Range("D29").Value = "'17.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.83%  "
$ws.Range("D30").Value = "'0.05666"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E31").Value = "  +2.55%  "
$ws.Range("D32").Value = "'3.650"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.30%  "
$ws.Range("D33").Value = "'3.503"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.05%  "
$ws.Range("E34").Value = "  +5.50%  "
$ws.Range("D35").Value = "'0.9736"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.05%  "
$ws.Range("D36").Value = "'2.842"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.66%  "
$ws.Range("D37").Value = "'2.436"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.35%  "
$ws.Range("D38").Value = "'0.5979"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.41%  "
$ws.Range("D39").Value = "'0.01669"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.83%  "
$ws.Range("D40").Value = "'5.908"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.33%  "
$ws.Range("E41").Value = "  +2.84%  "
$ws.Range("D42").Value = "'1.048.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.75%  "
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").Value = "'101.55"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").Value = "'1.872.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.11%  "
$ws.Range("E46").Value = "  +9.51%  "
$ws.Range("D47").Value = "'59.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.52%  "
$ws.Range("D48").Value = "'8.262"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.00%  "
$ws.Range("D49").Value = "'0.4428"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.94%  "
$ws.Range("D50").Value = "'0.05332"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("D51").Value = "'0.9995"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.64%  "
